# Add newly-collected runsheet rows for participants 46-50 (rows 47-51),
# then remove the now-unused "include/exclude" + "reason_excluded" columns (I:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 47 (SPEED_ACC_NOISE_46) ---
$ws.Range("B47").Value = 42964
$ws.Range("D47").Value = "F"
$ws.Range("E47").Value = 5
$ws.Range("H47").Value = "went well!"

# --- Row 48 (SPEED_ACC_NOISE_47) ---
$ws.Range("B48").Value = 42964
$ws.Range("D48").Value = "M"
$ws.Range("E48").Value = 4
$ws.Range("H48").Value = "saw me through the screen and got very distracted, eyetracker didn't pick up eyes at the end"

# --- Row 49 (SPEED_ACC_NOISE_48) ---
$ws.Range("B49").Value = 42964
$ws.Range("D49").Value = "M"
$ws.Range("E49").Value = 5
$ws.Range("H49").Value = "went well!"

# --- Row 50 (SPEED_ACC_NOISE_49) ---
$ws.Range("B50").Value = 42965
$ws.Range("D50").Value = "M"
$ws.Range("E50").Value = 4
$ws.Range("H50").Value = "went well!"

# --- Row 51 (SPEED_ACC_NOISE_50) ---
$ws.Range("B51").Value = 42965
$ws.Range("D51").Value = "F"
$ws.Range("E51").Value = 4
$ws.Range("H51").Value = "went well!"

# Remove the "exclude"/"reason_excluded" columns entirely - no longer tracked.
$ws.Range("I1:J66").Delete()

# Restore the usual view state (scrolled down a bit, B52 selected).
$ws.Activate()
$ws.Range("B52").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
